$p = $ppt.ActivePresentation

# Insert a new "Title and Content" slide at position 16 (pushes the
# existing "Next Steps" slide and everything after it down by one).
$s = $p.Slides.Add(16, 2)

# Leave the title placeholder empty (shape 1 - "Title 1").

# Fill the content placeholder (shape 2 - "Content Placeholder 2") with
# the pasted review text ("Google with legit reviews").
$s.Shapes.Item(2).TextFrame.TextRange.Text = '<span jstcache="142" aria-haspopup="true" role="button" tabindex="0" jsaction="pane.profile-stats.showStats;keydown:pane.profile-stats.showStats" class="uOKFHc-n1UuX-header-UjZuef uOKFHc-n1UuX-header-d6wfac-ibnC6b" jsan="7.uOKFHc-n1UuX-header-UjZuef,t-fUhMXXhk9tM,0.aria-haspopup,7.uOKFHc-n1UuX-header-d6wfac-ibnC6b,0.role,0.tabindex,0.jsaction">Local Guide · Level 8</span>'
